# "work to have multitimes and multiarc"
#
# 1. Rename the "Arcs" sheet to "Arc0".
# 2. Duplicate it as "Arc1" (same shape/styles), then flatten the per-arc
#    columns (C:H) of rows 11-18 to match column B (uniform across all arcs)
#    and bump the arc count (row 3) from 1 to 2.
# 3. Cameras!C6 (boresight elevation, camera 1) gets a refined value.
# 4. Arc1 becomes the active/selected sheet, matching the recorded view state.

$wb = $excel.ActiveWorkbook

# --- 1. Rename Arcs -> Arc0 -------------------------------------------------
$arc0 = $wb.Worksheets.Item("Arcs")
$arc0.Name = "Arc0"

# --- 2. Duplicate Arc0 -> Arc1, right after Arc0 ----------------------------
$arc0.Copy($null, $arc0)
$arc1 = $wb.Worksheets.Item($arc0.Index + 1)
$arc1.Name = "Arc1"

# Number of arcs (row 3, shared-string label "X0km") goes from 1 to 2 for
# every arc column.
$arc1.Range("B3:H3").Value2 = 2

# Flatten the per-arc-column values (C:H) onto column B's value for every
# per-arc numeric row, so all arcs share identical parameters.
foreach ($row in 11..18) {
    $bVal = $arc1.Cells.Item($row, 2).Value2
    $arc1.Range($arc1.Cells.Item($row, 3), $arc1.Cells.Item($row, 8)).Value2 = $bVal
}

# Cosmetic: Arc1 keeps no frozen header row/col and no custom tab colour in
# the recorded view; approximate that here.
$arc1.Activate()
$excel.ActiveWindow.FreezePanes = $false
$arc1.Range("C11:H18").Select()

# --- 3. Cameras!C6 updated value --------------------------------------------
$cameras = $wb.Worksheets.Item("Cameras")
$cameras.Range("C6").Value2 = 88.0172525718237

# --- 4. Final active sheet/selection state ----------------------------------
$arc1.Activate()
$arc1.Range("C11").Select()
